$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.324.65'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '2.477.91'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.08'
$ws.Range("E5").Value = '  -2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.60'
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.559'
$ws.Range("E8").Value = '  -1.32%  '
$ws.Range("D9").Value = '2.492.31'
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("E10").Value = '  -3.10%  '
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.340'
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("D14").Value = '2.918.10'
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").Value = '58.256.71'
$ws.Range("E15").Value = '  -0.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.18'
$ws.Range("E16").Value = '  -1.98%  '
$ws.Range("E17").Value = '  -1.88%  '
$ws.Range("D18").Value = '2.481.96'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("E19").Value = '  -2.71%  '
$ws.Range("E20").Value = '  -1.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.33'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.74'
$ws.Range("E23").Value = '  -3.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.52'
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.995'
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("E28").Value = '  -1.47%  '
$ws.Range("D29").Value = '0.0₃0752'
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.27'
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("E31").Value = '  -2.54%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.32'
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.19'
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.12'
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.61'
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("E40").Value = '  -2.85%  '
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("E42").Value = '  +3.54%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '275.17'
$ws.Range("E43").Value = '  -2.03%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.46'
$ws.Range("E44").Value = '  -2.83%  '
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '124.49'
$ws.Range("E46").Value = '  -3.93%  '
$ws.Range("E47").Value = '  -1.38%  '
$ws.Range("E48").Value = '  -1.52%  '
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.07'
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("D51").Value = '1.740.28'
$ws.Range("E51").Value = '  -0.77%  '
